$d = $word.ActiveDocument

# Paragraph that currently only holds the _GoBack bookmark.
$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Add a new, empty paragraph right after it (before the final sectPr).
$bookmarkPara.Range.InsertParagraphAfter()

# Turn the (now second-to-last) bookmark paragraph into a "List Paragraph"
# styled paragraph containing the text "List Paragraph", keeping the
# bookmark start/end at the end of the paragraph.
$bookmarkPara.Range.InsertBefore("List Paragraph")
$bookmarkPara.Style = "List Paragraph"

# Match the style definition minted by Word for List Paragraph.
$listStyle = $d.Styles.Item("List Paragraph")
$listStyle.Priority = 34
$listStyle.NoSpaceBetweenParagraphsOfSameStyle = $true
$listStyle.ParagraphFormat.LeftIndent = 36
